$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 605.3333
$ws.Range("J43").Value = 678
$ws.Range("L43").Value = 678
$ws.Range("N43").Value = -816
$ws.Range("H62").Value = 4181.0625
$ws.Range("I62").Value = 4181.0625
$ws.Range("K62").Value = 4181.0625
$ws.Range("M62").Value = -3557.0625
$ws.Range("H65").Value = 4181.0625
$ws.Range("I65").Value = 4181.0625
$ws.Range("K65").Value = 20905.3125
$ws.Range("M65").Value = -17785.3125
$ws.Range("H137").Value = 3540.5
$ws.Range("I137").Value = 4310.75
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 12932.25
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -10382.25
$ws.Range("N137").Value = -11100
$ws.Range("H138").Value = 1923
$ws.Range("I138").Value = 1524.8096
$ws.Range("J138").Value = 3316.6667
$ws.Range("K138").Value = 4574.4288
$ws.Range("L138").Value = 9950.000100000001
$ws.Range("M138").Value = 565.5712000000003
$ws.Range("N138").Value = -20230.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -50976
$ws.Range("H55").Value = 35000
$ws.Range("I55").Value = 30000
$ws.Range("J55").Value = 40000
$ws.Range("K55").Value = 30000
$ws.Range("L55").Value = 40000
$ws.Range("M55").Value = -29685
$ws.Range("N55").Value = -40630
$ws.Range("H61").Value = 3291.6843
$ws.Range("I61").Value = 1782.125
$ws.Range("J61").Value = 11342.667
$ws.Range("K61").Value = 1782.125
$ws.Range("L61").Value = 11342.667
$ws.Range("M61").Value = -1570.125
$ws.Range("N61").Value = -11766.667
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680
$ws.Range("H122").Value = 2260.0688
$ws.Range("I122").Value = 1627.7916
$ws.Range("J122").Value = 5295
$ws.Range("K122").Value = 4883.3748
$ws.Range("L122").Value = 15885
$ws.Range("M122").Value = -2433.3748
$ws.Range("N122").Value = -20785
$ws.Range("H136").Value = 3291.6843
$ws.Range("I136").Value = 1782.125
$ws.Range("J136").Value = 11342.667
$ws.Range("K136").Value = 5346.375
$ws.Range("L136").Value = 34028.001
$ws.Range("M136").Value = -2796.375
$ws.Range("N136").Value = -39128.001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2014
$ws.Range("I86").Value = 1669.2858
$ws.Range("J86").Value = 2416.1667
$ws.Range("K86").Value = 1669.2858
$ws.Range("L86").Value = 2416.1667
$ws.Range("M86").Value = -546.2858000000001
$ws.Range("N86").Value = -4662.1667
$ws.Range("H89").Value = 2014
$ws.Range("I89").Value = 1669.2858
$ws.Range("J89").Value = 2416.1667
$ws.Range("K89").Value = 8346.429
$ws.Range("L89").Value = 12080.8335
$ws.Range("M89").Value = -2730.429
$ws.Range("N89").Value = -23312.8335
$ws.Range("H107").Value = 2220
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 2220
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2220
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -6060

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3166.625
$ws.Range("I16").Value = 2619
$ws.Range("J16").Value = 7000
$ws.Range("K16").Value = 2619
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = -2332
$ws.Range("N16").Value = -7574
$ws.Range("H62").Value = 11541.615
$ws.Range("I62").Value = 2604.1
$ws.Range("J62").Value = 41333.332
$ws.Range("K62").Value = 2604.1
$ws.Range("L62").Value = 41333.332
$ws.Range("M62").Value = -1980.1
$ws.Range("N62").Value = -42581.332
$ws.Range("H65").Value = 11541.615
$ws.Range("I65").Value = 2604.1
$ws.Range("J65").Value = 41333.332
$ws.Range("K65").Value = 13020.5
$ws.Range("L65").Value = 206666.66
$ws.Range("M65").Value = -9900.5
$ws.Range("N65").Value = -212906.66
$ws.Range("H99").Value = 3061.7827
$ws.Range("I99").Value = 2528.1333
$ws.Range("K99").Value = 2528.1333
$ws.Range("M99").Value = -1030.1333
$ws.Range("H109").Value = 27000.8
$ws.Range("J109").Value = 27000.8
$ws.Range("L109").Value = 27000.8
$ws.Range("N109").Value = -29080.8
$ws.Range("H113").Value = 3166.625
$ws.Range("I113").Value = 2619
$ws.Range("J113").Value = 7000
$ws.Range("K113").Value = 2619
$ws.Range("L113").Value = 7000
$ws.Range("M113").Value = -449
$ws.Range("N113").Value = -11340
$ws.Range("H122").Value = 41669040
$ws.Range("I122").Value = 62500910
$ws.Range("K122").Value = 187502730
$ws.Range("M122").Value = -187500280
$ws.Range("H126").Value = 3061.7827
$ws.Range("I126").Value = 2528.1333
$ws.Range("K126").Value = 7584.3999
$ws.Range("M126").Value = -5114.3999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 20000000
$ws.Range("I24").Value = 20000000
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 20000000
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -19999827
$ws.Range("N24").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H29").Value = 10266667
$ws.Range("I29").Value = 15000000
$ws.Range("J29").Value = 800000
$ws.Range("K29").Value = 15000000
$ws.Range("L29").Value = 800000
$ws.Range("N29").Value = -800580
$ws.Range("M29").Value = -14999710
$ws.Range("H107").Value = 689.619
$ws.Range("I107").Value = 634.1818
$ws.Range("J107").Value = 750.6
$ws.Range("K107").Value = 634.1818
$ws.Range("L107").Value = 750.6
$ws.Range("M107").Value = 1285.8182
$ws.Range("N107").Value = -4590.6
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = 0
$ws.Range("H132").Value = 3257.6365
$ws.Range("I132").Value = 2870.6667
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 8612.000100000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -6082.000100000001
$ws.Range("N132").Value = -20057

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2090.3333
$ws.Range("I7").Value = 1834.6666
$ws.Range("J7").Value = 2601.6667
$ws.Range("K7").Value = 1834.6666
$ws.Range("L7").Value = 2601.6667
$ws.Range("M7").Value = -1722.6666
$ws.Range("N7").Value = -2825.6667
$ws.Range("H22").Value = 585.087
$ws.Range("I22").Value = 231
$ws.Range("J22").Value = 971.36365
$ws.Range("K22").Value = 231
$ws.Range("L22").Value = 971.36365
$ws.Range("M22").Value = 64
$ws.Range("N22").Value = -1561.36365
$ws.Range("H27").Value = 585.087
$ws.Range("I27").Value = 231
$ws.Range("J27").Value = 971.36365
$ws.Range("K27").Value = 231
$ws.Range("L27").Value = 971.36365
$ws.Range("M27").Value = -124
$ws.Range("N27").Value = -1185.36365
$ws.Range("H126").Value = 2090.3333
$ws.Range("I126").Value = 1834.6666
$ws.Range("J126").Value = 2601.6667
$ws.Range("K126").Value = 5503.9998
$ws.Range("L126").Value = 7805.000100000001
$ws.Range("M126").Value = -3033.9998
$ws.Range("N126").Value = -12745.0001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3491.158
$ws.Range("I62").Value = 2626.2307
$ws.Range("J62").Value = 5365.1665
$ws.Range("K62").Value = 2626.2307
$ws.Range("L62").Value = 5365.1665
$ws.Range("M62").Value = -2002.2307
$ws.Range("N62").Value = -6613.1665
$ws.Range("H65").Value = 3491.158
$ws.Range("I65").Value = 2626.2307
$ws.Range("J65").Value = 5365.1665
$ws.Range("K65").Value = 13131.1535
$ws.Range("L65").Value = 26825.8325
$ws.Range("M65").Value = -10011.1535
$ws.Range("N65").Value = -33065.8325
$ws.Range("H109").Value = 24597.5
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 24597.5
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 24597.5
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -27371.5
$ws.Range("H132").Value = 2495.742
$ws.Range("I132").Value = 1960.3125
$ws.Range("J132").Value = 3066.8667
$ws.Range("K132").Value = 5880.9375
$ws.Range("L132").Value = 9200.6001
$ws.Range("M132").Value = -3350.9375
$ws.Range("N132").Value = -14260.6001
